$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 53; this shifts the existing rows 53-67
# down to 54-68 (and the sheet dimension grows from R67 to R68).
$ws.Rows("53").Insert()

# Populate the newly inserted row 53 with the new record.
$ws.Range("A53").Value = 5
$ws.Range("B53").Value = "Macroferia Regional de Talca"
$ws.Range("C53").Value = "Maule"
$ws.Range("D53").Value = 44837
$ws.Range("E53").Value = 7
$ws.Range("F53").Value = 300000000
$ws.Range("G53").Value = "Espárragos"
$ws.Range("H53").Value = "Sin especificar"
$ws.Range("I53").Value = "Primera"
$ws.Range("J53").Value = 2000
$ws.Range("K53").Value = 1600
$ws.Range("L53").Value = 1600
$ws.Range("M53").Value = 1600
$ws.Range("N53").Value = '$/kilo'
$ws.Range("O53").Value = "Provincia de Linares"
$ws.Range("P53").Value = 1600
$ws.Range("Q53").Value = 1
$ws.Range("R53").Value = "Hortaliza"
